# B6-PowerPoint.pptx edit
#
# 1. Three tables (on slides 14, 15 and 16) switch from the locally-defined
#    "Table_0" table style ({64AB2897-85DB-4337-90BE-549DBE76E6F0}) to the
#    built-in "Medium Style 2" table style
#    ({45686A7C-952D-4BE4-B6A4-CDF838DA80D6}).
#
# 2. The deck's theme palette is swapped: the slide master's theme
#    ("Integral" / "Red Violet") is replaced with the plain default
#    "Office Theme" palette (which, before the edit, only the notes master
#    used). We recreate that by writing the 12 target colours into the
#    slide theme's colour scheme.

$p = $ppt.ActivePresentation

# --- 1) Retarget the three tables to the built-in table style ----------
$newTableStyleId = "{45686A7C-952D-4BE4-B6A4-CDF838DA80D6}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2) Swap the slide theme's colour scheme to the "Office" palette ---
# Index order matches the OOXML <a:clrScheme> child order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# RGB values below use the COM RGB() packing (R + G*256 + B*65536), derived
# from the target srgbClr hex values 000000/FFFFFF/44546A/E7E6E6/5B9BD5/
# ED7D31/A5A5A5/FFC000/4472C4/70AD47/0563C1/954F72.
$officeColors = @{
    1  = 0        # dk1      000000
    2  = 16777215 # lt1      FFFFFF
    3  = 6968388  # dk2      44546A
    4  = 15132391 # lt2      E7E6E6
    5  = 13998939 # accent1  5B9BD5
    6  = 3243501  # accent2  ED7D31
    7  = 10855845 # accent3  A5A5A5
    8  = 49407    # accent4  FFC000
    9  = 12874308 # accent5  4472C4
    10 = 4697456  # accent6  70AD47
    11 = 12673797 # hlink    0563C1
    12 = 7491477  # folHlink 954F72
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($idx = 1; $idx -le 12; $idx++) {
    $themeColors.Item($idx).RGB = $officeColors[$idx]
}
